$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.485378203517801
$ws.Range("I2").Value = 0.437672734277873
$ws.Range("K2").Value = 0.36062834418311
$ws.Range("L2").Value = 0.647351101749226
$ws.Range("N2").Value = 0.450957959509906

$ws.Range("B3").Value = 0.372303915214208
$ws.Range("J3").Value = 0.435665506975429
$ws.Range("K3").Value = 0.357209918266145
$ws.Range("L3").Value = 0.230868411160313
$ws.Range("N3").Value = 0.366083392221669

$ws.Range("B4").Value = 0.368431562488216
$ws.Range("K4").Value = 0.268475855361453
$ws.Range("L4").Value = 0.34323165564695
$ws.Range("N4").Value = 0.354982855614853

$ws.Range("B5").Value = 0.336207928322194
$ws.Range("K5").Value = 0.240441928542887
$ws.Range("L5").Value = 0.536926339855414
$ws.Range("N5").Value = 0.308524690287428

$ws.Range("B6").Value = 0.28887016021309
$ws.Range("K6").Value = 0.221115870507474
$ws.Range("L6").Value = 0.17381964262529
$ws.Range("N6").Value = 0.291228468936834

$ws.Range("B7").Value = 0.220446803986899
$ws.Range("K7").Value = 0.183257456612112
$ws.Range("L7").Value = 0.19554617999603
$ws.Range("N7").Value = 0.217448055217537

$ws.Range("B8").Value = 0.184742937735908
$ws.Range("K8").Value = 0.141433314757512
$ws.Range("L8").Value = 0.185034790211497
$ws.Range("N8").Value = 0.178603799379371
